# Product_Backlog.xlsx edit: add a new backlog item to the "Sprint 6"
# group (row 15) for "view invoice list" user story, push the existing
# "view statistics" item down, and leave a fresh blank row at the bottom
# of the tracker (row 24). Also nudge the view state (scrolled to A10,
# selection on E15) the way it looks right after the edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 15. Excel automatically
# shifts the row contents/heights/styles of rows 15-23 down to 16-24,
# extends the A13:A15 / A16:A18 merged groups to A13:A16 / A17:A19, and
# grows the sheet dimension to A3:H24 for us.
$ws.Rows("15").Insert()

# Row 15 is the new blank row; fill in the new backlog item. It belongs
# to the same "Sprint 6" merged group as rows 13-16, so column A stays
# blank, and it is the 3rd sub-item of that sprint (the old 3rd item,
# "view statistics", becomes the 4th item down on row 16).
$ws.Range("A15").Value = $null
$ws.Range("B15").Value = 3
$ws.Range("C15").Value = "Là một khách hàng, tôi muốn xem lại danh sách các hóa đơn đã đặt"
$ws.Range("D15").Value = "Xây dưng trang quản lý hóa đơn của khách hàng"
$ws.Rows("15").RowHeight = 43.5

# The item that used to be row 15 ("view statistics") is now row 16 and
# becomes sub-item 4 of the Sprint-6 group; its 58pt row height already
# carried down with the rest of the row's formatting when row 15 was
# inserted above it, so only the number needs updating.
$ws.Range("B16").Value = 4

# View state: scrolled down a bit and selection sitting on E15 (instead
# of F15) after the edit.
$ws.Range("E15").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 10
